$wb = $excel.ActiveWorkbook

# The same update (refreshed "想去人数" / interested-count figures) needs to be
# applied to both the "展览" sheet and the "全部类型" sheet, which mirror the
# same underlying data.
$sheetNames = @("展览", "全部类型")

# Map of cell address -> new value for column F ("想去人数")
$updates = @{
    "F2"  = 692
    "F3"  = 5
    "F9"  = 3399
    "F10" = 4286
    "F12" = 129
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellAddr in $updates.Keys) {
        $ws.Range($cellAddr).Value = $updates[$cellAddr]
    }
}
